$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lương")

# Delete the "tại HỆ THỐNG" block (rows 4-10) and the "Tổng lương tại HỆ THỐNG"
# row (originally row 35) — these are the discount-ratio summary rows for the
# "HỆ THỐNG" location that the commit removes.
# Delete bottom-up so earlier row numbers stay valid while later ones are removed.
$ws.Range("A35").EntireRow.Delete()
$ws.Range("A4:A10").EntireRow.Delete()

# After the deletion, the remaining rows have shifted up:
#  old row 2 (Ngày công)                -> still row 2
#  old row 3 (Phụ cấp)                  -> still row 3
#  old row 11 (Lương cơ bản tại CẦN THƠ)-> now row 4
#  old row 36 (Tổng lương tại CẦN THƠ)  -> now row 28
#  old row 39 (Tổng lương)              -> now row 31
$ws.Range("B2").Value = 25
$ws.Range("B3").Value = 875000
$ws.Range("B4").Value = 5357142.857142857
$ws.Range("B28").Value = 6232142.857142857
$ws.Range("B31").Value = 6232142.857142857
